# daily auto push: 2026-01-12 13:51 UTC
# Insert a new log row for 2026/01/12 (月) at row 622, pushing the
# existing rows 622:663 down to 623:664.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 622 (old rows 622-663
# shift down to 623-664).
$ws.Rows.Item(622).Insert()

# Seed the new row by copying the row directly above it (row 621 is
# already "2026/01/12" / "月" with the same text-typed date column),
# which keeps A/B as plain text cells (no auto date conversion, no new
# number-format styles) instead of assigning literal strings.
$ws.Range("A621:D621").Copy()
$ws.Range("A622:D622").PasteSpecial()
$excel.CutCopyMode = $false

# Only the hour/time value differs for the new entry.
$ws.Range("C622").Value = 19
